$d = $word.ActiveDocument
$d.Content.Find.Execute("70×12=840", $true, $false, $false, $false, $false, $true, 1, $false, "80×26=2080", 2) | Out-Null
$d.Content.Find.Execute("72×55=3960", $true, $false, $false, $false, $false, $true, 1, $false, "75×43=3225", 2) | Out-Null
$d.Content.Find.Execute("14×62=868", $true, $false, $false, $false, $false, $true, 1, $false, "87×55=4785", 2) | Out-Null
$d.Content.Find.Execute("48×21=1008", $true, $false, $false, $false, $false, $true, 1, $false, "56×14=784", 2) | Out-Null
$d.Content.Find.Execute("46×15=690", $true, $false, $false, $false, $false, $true, 1, $false, "33×60=1980", 2) | Out-Null
$d.Content.Find.Execute("87×94=8178", $true, $false, $false, $false, $false, $true, 1, $false, "75×14=1050", 2) | Out-Null
$d.Content.Find.Execute("29×15=435", $true, $false, $false, $false, $false, $true, 1, $false, "88×23=2024", 2) | Out-Null
$d.Content.Find.Execute("85×77=6545", $true, $false, $false, $false, $false, $true, 1, $false, "75×28=2100", 2) | Out-Null
$d.Content.Find.Execute("87×69=6003", $true, $false, $false, $false, $false, $true, 1, $false, "24×62=1488", 2) | Out-Null
$d.Content.Find.Execute("74×57=4218", $true, $false, $false, $false, $false, $true, 1, $false, "13×14=182", 2) | Out-Null
$d.Content.Find.Execute("16×59=944", $true, $false, $false, $false, $false, $true, 1, $false, "18×58=1044", 2) | Out-Null
$d.Content.Find.Execute("13×26=338", $true, $false, $false, $false, $false, $true, 1, $false, "54×54=2916", 2) | Out-Null
$d.Content.Find.Execute("64×96=6144", $true, $false, $false, $false, $false, $true, 1, $false, "67×38=2546", 2) | Out-Null
$d.Content.Find.Execute("47×18=846", $true, $false, $false, $false, $false, $true, 1, $false, "87×13=1131", 2) | Out-Null
$d.Content.Find.Execute("11×49=539", $true, $false, $false, $false, $false, $true, 1, $false, "39×69=2691", 2) | Out-Null
$d.Content.Find.Execute("84×90=7560", $true, $false, $false, $false, $false, $true, 1, $false, "66×48=3168", 2) | Out-Null
$d.Content.Find.Execute("62×84=5208", $true, $false, $false, $false, $false, $true, 1, $false, "83×52=4316", 2) | Out-Null
$d.Content.Find.Execute("12×42=504", $true, $false, $false, $false, $false, $true, 1, $false, "88×55=4840", 2) | Out-Null
$d.Content.Find.Execute("21×78=1638", $true, $false, $false, $false, $false, $true, 1, $false, "38×57=2166", 2) | Out-Null
$d.Content.Find.Execute("24×50=1200", $true, $false, $false, $false, $false, $true, 1, $false, "60×57=3420", 2) | Out-Null
$d.Content.Find.Execute("50×61=3050", $true, $false, $false, $false, $false, $true, 1, $false, "41×50=2050", 2) | Out-Null
$d.Content.Find.Execute("36×84=3024", $true, $false, $false, $false, $false, $true, 1, $false, "99×43=4257", 2) | Out-Null
$d.Content.Find.Execute("94×95=8930", $true, $false, $false, $false, $false, $true, 1, $false, "51×42=2142", 2) | Out-Null
$d.Content.Find.Execute("16×30=480", $true, $false, $false, $false, $false, $true, 1, $false, "47×71=3337", 2) | Out-Null
$d.Content.Find.Execute("74×58=4292", $true, $false, $false, $false, $false, $true, 1, $false, "18×95=1710", 2) | Out-Null
